$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated Neo4j "FilesTab" query (B4): drop the `file_type`/`File Type` and
# `demo.breed`/`Breed` columns per the corrected ICDC Breed 1-14 scripts.
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['American Staffordshire Terrier']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newQuery

# Row 4 shrinks now that the query text has two fewer lines.
$ws.Rows.Item(4).RowHeight = 217.5

# Scroll so row 4 is the window's top row, and move the selection from A4 to B4.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("B4").Select()
